$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SCA_N)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.677167643176708

# Row 3 (EA_N)
$ws.Range("B3").Value = 0.7141847117164766
$ws.Range("C3").Value = -0.6656932886364156
$ws.Range("D3").Value = -0.7613811385274181

# Row 4 (ENSO-mei_N)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = -0.6745126611206681
$ws.Range("D4").Value = -0.6786028277126004

# Row 5 (NAO_N)
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.63840281278613

# Row 6 (SCA_P)
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = -0.8038404870160419

# Row 7 (EA_P)
$ws.Range("B7").Value = -0.7493204663214563
$ws.Range("C7").Value = -0.6580824943572741
$ws.Range("D7").Value = 0.8381110819000569

# Row 8 (ENSO-mei_P)
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = -0.6300590212090541
$ws.Range("D8").Value = -0.8343932680828058

# Row 9 (NAO_P)
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = -0.6921262015650442
